$wb = $excel.ActiveWorkbook

# --- Dollar sheet: insert a new "Invoiced Month" column before InvoiceNumber ---
$wsDollar = $wb.Worksheets.Item("Dollar")
$wsDollar.Columns.Item(2).Insert()
$wsDollar.Range("B1").Value = "Invoiced Month"

# Update the "Dollars" defined name so it covers the new column range
$wb.Names.Item("Dollars").RefersTo = "=Dollar!`$A`$1:`$Q`$1"

# --- Clients sheet: drop the sample rows, relabel the remaining header ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Rows.Item(2).Resize(3).Delete()
$wsClients.Range("A1").Value = "Client"

# Highlight duplicate client names with the standard "duplicate values" style
$range = $wsClients.Range("A2:XFD242")
$fc = $range.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$fc.Priority = 2
